# Update the document date and the division problems/answers in the table.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-06-19 Thursday"; new = "2025-06-20 Friday"},
    @{old = "295÷7=42, 1";  new = "680÷9=75, 5"},
    @{old = "221÷4=55, 1";  new = "233÷3=77, 2"},
    @{old = "354÷4=88, 2";  new = "305÷4=76, 1"},
    @{old = "951÷2=475, 1"; new = "867÷4=216, 3"},
    @{old = "425÷2=212, 1"; new = "473÷7=67, 4"},
    @{old = "935÷9=103, 8"; new = "610÷2=305, 0"},
    @{old = "978÷9=108, 6"; new = "763÷2=381, 1"},
    @{old = "976÷4=244, 0"; new = "402÷7=57, 3"},
    @{old = "237÷7=33, 6";  new = "522÷7=74, 4"},
    @{old = "686÷8=85, 6";  new = "978÷8=122, 2"},
    @{old = "985÷8=123, 1"; new = "270÷6=45, 0"},
    @{old = "200÷6=33, 2";  new = "173÷8=21, 5"},
    @{old = "643÷2=321, 1"; new = "389÷6=64, 5"},
    @{old = "841÷7=120, 1"; new = "311÷2=155, 1"},
    @{old = "673÷2=336, 1"; new = "624÷7=89, 1"},
    @{old = "524÷6=87, 2";  new = "623÷5=124, 3"},
    @{old = "876÷5=175, 1"; new = "427÷7=61, 0"},
    @{old = "850÷4=212, 2"; new = "513÷4=128, 1"},
    @{old = "838÷9=93, 1";  new = "346÷7=49, 3"},
    @{old = "492÷4=123, 0"; new = "577÷5=115, 2"},
    @{old = "923÷6=153, 5"; new = "127÷2=63, 1"},
    @{old = "912÷8=114, 0"; new = "169÷4=42, 1"},
    @{old = "842÷2=421, 0"; new = "781÷7=111, 4"},
    @{old = "940÷6=156, 4"; new = "762÷3=254, 0"},
    @{old = "305÷3=101, 2"; new = "673÷6=112, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
